# Apply cryptos list update (values refreshed by scraper run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.064.60"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.792.78"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "222.40"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.550"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.79"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -3.51%  "
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0716"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +4.61%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.052.08"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.794.89"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.74"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.059.31"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.23"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.07"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.59"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.08"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.49"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.69"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.51"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.410.03"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.643"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.945"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +4.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.04"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.95"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0493"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.950.77"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("E47").Value = "  -3.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.62"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.84"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0121"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.09%  "
